$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.453645056076748
$ws.Range("L2").Value = 0.498518803138412

$ws.Range("B3").Value = 0.318030317882592
$ws.Range("L3").Value = 0.22699232463484

$ws.Range("B4").Value = 0.172370897143246
$ws.Range("E4").Value = 0.166158467826809
$ws.Range("L4").Value = 0.195595011071219
